# Swap the entire content of row 2 and row 3 (all columns A-E), including
# the hyperlink attached to column E, so that the "Alleged Darwin shooter
# Benjamin Hoffmann..." article becomes row 2 and the "Chilling CCTV..."
# article becomes row 3. (Rows 4 and 5 are untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current (pre-edit) row values for all 4 data rows / all columns ---
$a2 = $ws.Cells.Item(2, 1).Value()
$b2 = $ws.Cells.Item(2, 2).Value()
$c2 = $ws.Cells.Item(2, 3).Value()
$d2 = $ws.Cells.Item(2, 4).Value()
$e2 = $ws.Cells.Item(2, 5).Value()

$a3 = $ws.Cells.Item(3, 1).Value()
$b3 = $ws.Cells.Item(3, 2).Value()
$c3 = $ws.Cells.Item(3, 3).Value()
$d3 = $ws.Cells.Item(3, 4).Value()
$e3 = $ws.Cells.Item(3, 5).Value()

# the uri column (E) text *is* the hyperlink target address for these rows
$e4 = $ws.Cells.Item(4, 5).Value()
$e5 = $ws.Cells.Item(5, 5).Value()

# --- write row 2 <- old row 3, and row 3 <- old row 2 (columns A-E) ---
$ws.Cells.Item(2, 1).Value = $a3
$ws.Cells.Item(2, 2).Value = $b3
$ws.Cells.Item(2, 3).Value = $c3
$ws.Cells.Item(2, 4).Value = $d3
$ws.Cells.Item(2, 5).Value = $e3

$ws.Cells.Item(3, 1).Value = $a2
$ws.Cells.Item(3, 2).Value = $b2
$ws.Cells.Item(3, 3).Value = $c2
$ws.Cells.Item(3, 4).Value = $d2
$ws.Cells.Item(3, 5).Value = $e2

# --- rebuild the hyperlinks collection so E2/E3 point at the swapped URLs
#     (this engine only supports clearing + re-adding the whole collection) ---
$hl = $ws.Hyperlinks
$hl.Delete()

$hl.Add($ws.Range("E2"), $e3)
$hl.Add($ws.Range("E3"), $e2)
$hl.Add($ws.Range("E4"), $e4)
$hl.Add($ws.Range("E5"), $e5)

# restore the Hyperlink cell style on the uri column
$ws.Range("E2:E5").Style = "Hyperlink"
